$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Ensures the value is stored as text even when it looks like a number,
    # while keeping the cell's original (unstyled) appearance.
    $ws.Range($cellRef).Value = "'" + $value
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "59.726.83"
$ws.Range("E2").Value = "  +3.16%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.417.13"
$ws.Range("E3").Value = "  +2.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "551.55"
$ws.Range("E5").Value = "  +0.40%  "

# Row 6 - Solana
Set-TextValue "D6" "136.95"
$ws.Range("E6").Value = "  +3.52%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - XRP
Set-TextValue "D8" "0.589"
$ws.Range("E8").Value = "  +3.73%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  +0.08%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  +0.70%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  -2.24%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.16%  "

# Row 13 - Avalanche
$ws.Range("E13").Value = "  +3.80%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "2.849.24"
$ws.Range("E14").Value = "  +2.87%  "

# Row 15 - WrappedBTC
Set-TextValue "D15" "59.688.79"
$ws.Range("E15").Value = "  +3.34%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.01%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.412.05"
$ws.Range("E17").Value = "  +3.52%  "

# Row 18 - Chainlink
$ws.Range("E18").Value = "  +2.63%  "

# Row 19 - Polkadot
$ws.Range("E19").Value = "  +2.10%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "330.45"
$ws.Range("E20").Value = "  +0.21%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -2.90%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.05%  "

# Row 23 - Litecoin
Set-TextValue "D23" "66.00"
$ws.Range("E23").Value = "  +3.38%  "

# Row 24 - Kaspa
$ws.Range("E24").Value = "  +1.84%  "

# Row 25 - InternetComputer(DFINITY)
Set-TextValue "D25" "8.72"
$ws.Range("E25").Value = "  +5.74%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.10%  "

# Row 27 - Fetch.AI
Set-TextValue "D27" "1.36"
$ws.Range("E27").Value = "  +3.60%  "

# Row 28 - PEPE
Set-TextValue "D28" "0.0₃0774"
$ws.Range("E28").Value = "  +5.18%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +0.27%  "

# Row 30 - Monero
Set-TextValue "D30" "170.20"
$ws.Range("E30").Value = "  -0.49%  "

# Row 31 - Aptos
Set-TextValue "D31" "6.13"
$ws.Range("E31").Value = "  +0.13%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "18.65"
$ws.Range("E32").Value = "  +1.70%  "

# Row 35 - ImmutableX
$ws.Range("E35").Value = "  +4.32%  "

# Row 36 - FirstDigitalUSD
$ws.Range("E36").Value = "  +0.06%  "

# Row 37 - NEARProtocol
$ws.Range("E37").Value = "  +1.47%  "

# Row 38 - Stacks
$ws.Range("E38").Value = "  +0.45%  "

# Row 39 - OKB
Set-TextValue "D39" "39.55"
$ws.Range("E39").Value = "  -2.00%  "

# Row 40 - PolygonEcosystemToken
Set-TextValue "D40" "0.411"
$ws.Range("E40").Value = "  -5.98%  "

# Row 41 - Bittensor
Set-TextValue "D41" "313.39"

# Row 42 - Filecoin
$ws.Range("E42").Value = "  +0.63%  "

# Row 43 - Aave
Set-TextValue "D43" "138.90"
$ws.Range("E43").Value = "  -2.23%  "

# Row 44 - Stellar
Set-TextValue "D44" "0.0970"
$ws.Range("E44").Value = "  +1.78%  "

# Row 45 - Hedera
$ws.Range("E45").Value = "  +0.99%  "

# Row 46 - InjectiveProtocol
$ws.Range("E46").Value = "  +4.26%  "

# Row 47 - Mantle
Set-TextValue "D47" "0.579"
$ws.Range("E47").Value = "  +2.23%  "

# Row 48 - VeChain
$ws.Range("E48").Value = "  +1.25%  "

# Row 49 - Polygon
$ws.Range("E49").Value = "  -9.12%  "

# Row 50 - EnergySwap
Set-TextValue "D50" "17.58"
$ws.Range("E50").Value = "  +0.93%  "

# Row 51 - WhiteBITCoin
Set-TextValue "D51" "11.06"
$ws.Range("E51").Value = "  -0.05%  "
